# Applies the "Saldo_guide" update:
#  - Reference date moves from 2024-11-01 (serial 45597) to 2024-11-04 (serial 45600)
#    for every data row (column G, rows 2..274).
#  - A handful of balance rows get corrected Vl. Projetado / Saldo Previsto / Vl. Total
#    figures (columns D/E/H).
#  - The sheet tab name carries the new extraction timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the "Dt. Referencia" column (G) from 45597 to 45600 for all data rows.
$ws.Range("G2:G274").Value = 45600

# 2) Row 51 - Saldo Previsto / Vl. Total corrected from 11765.61 to 765.61
$ws.Range("E51").Value = 765.61
$ws.Range("H51").Value = 765.61

# 3) Row 107 - Saldo Previsto / Vl. Total corrected from 22323.97 to 323.97000000000003
$ws.Range("E107").Value = 323.97000000000003
$ws.Range("H107").Value = 323.97000000000003

# 4) Row 143 - Vl. Projetado, Saldo Previsto and Vl. Total all corrected
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 151.38999999999999
$ws.Range("H143").Value = 151.38999999999999

# 5) Row 230 - Saldo Previsto / Vl. Total corrected from 6397.75 to 397.75
$ws.Range("E230").Value = 397.75
$ws.Range("H230").Value = 397.75

# 6) Row 231 - Saldo Previsto / Vl. Total corrected from 673.78 to 602.96
$ws.Range("E231").Value = 602.96
$ws.Range("H231").Value = 602.96

# 7) Row 232 - Saldo Previsto / Vl. Total corrected from 33881.82 to 22881.82
$ws.Range("E232").Value = 22881.82
$ws.Range("H232").Value = 22881.82

# 8) Rename the sheet tab to reflect the new export timestamp.
$ws.Name = "IClientBalance-20241104-090629-"
